# Scheduled-runner refresh: update market/profit columns (H-N) on the
# Kujata_Profits leve-crafting sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values below are the freshly-pulled currentAveragePrice* / LevePrice* /
# LeveProfit* figures; a handful of rows also gain/lose their HQ-profit
# cell (N) or NQ-profit cell (M) depending on whether an HQ/NQ price now
# exists for that item.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 159.33333
$ws.Cells.Item(33, 9).Value = 106.07692
$ws.Cells.Item(33, 11).Value = 106.07692
$ws.Cells.Item(33, 13).Value = 122.92308
$ws.Cells.Item(40, 8).Value = 3599.875
$ws.Cells.Item(40, 9).Value = 3542.7144
$ws.Cells.Item(40, 10).Value = 4000
$ws.Cells.Item(40, 11).Value = 3542.7144
$ws.Cells.Item(40, 12).Value = 4000
$ws.Cells.Item(40, 13).Value = -3367.7144
$ws.Cells.Item(40, 14).Value = -4350
$ws.Cells.Item(86, 8).Value = 5042.5
$ws.Cells.Item(86, 9).Value = 3556.6667
$ws.Cells.Item(86, 10).Value = 9500
$ws.Cells.Item(86, 11).Value = 3556.6667
$ws.Cells.Item(86, 12).Value = 9500
$ws.Cells.Item(86, 13).Value = -2433.6667
$ws.Cells.Item(86, 14).Value = -11746
$ws.Cells.Item(89, 8).Value = 5042.5
$ws.Cells.Item(89, 9).Value = 3556.6667
$ws.Cells.Item(89, 10).Value = 9500
$ws.Cells.Item(89, 11).Value = 17783.3335
$ws.Cells.Item(89, 12).Value = 47500
$ws.Cells.Item(89, 13).Value = -12167.3335
$ws.Cells.Item(89, 14).Value = -58732
$ws.Cells.Item(98, 8).Value = 1533.7
$ws.Cells.Item(98, 9).Value = 1560.875
$ws.Cells.Item(98, 10).Value = 1425
$ws.Cells.Item(98, 11).Value = 1560.875
$ws.Cells.Item(98, 12).Value = 1425
$ws.Cells.Item(98, 13).Value = -62.875
$ws.Cells.Item(98, 14).Value = -4421
$ws.Cells.Item(112, 8).Value = 3702.5334
$ws.Cells.Item(112, 10).Value = 4103
$ws.Cells.Item(112, 12).Value = 12309
$ws.Cells.Item(112, 14).Value = -14525
$ws.Cells.Item(122, 8).Value = 1533.7
$ws.Cells.Item(122, 9).Value = 1560.875
$ws.Cells.Item(122, 10).Value = 1425
$ws.Cells.Item(122, 11).Value = 4682.625
$ws.Cells.Item(122, 12).Value = 4275
$ws.Cells.Item(122, 13).Value = -2232.625
$ws.Cells.Item(122, 14).Value = -9175
$ws.Cells.Item(129, 8).Value = 848.5161000000001
$ws.Cells.Item(129, 10).Value = 884.0714
$ws.Cells.Item(129, 12).Value = 2652.2142
$ws.Cells.Item(129, 14).Value = -12652.2142
$ws.Cells.Item(138, 8).Value = 2061.78
$ws.Cells.Item(138, 9).Value = 895.0833
$ws.Cells.Item(138, 10).Value = 2220.875
$ws.Cells.Item(138, 11).Value = 2685.2499
$ws.Cells.Item(138, 12).Value = 6662.625
$ws.Cells.Item(138, 13).Value = 2454.7501
$ws.Cells.Item(138, 14).Value = -16942.625
$ws.Cells.Item(141, 8).Value = 8056.8
$ws.Cells.Item(141, 9).Value = 9637.75
$ws.Cells.Item(141, 10).Value = 1733
$ws.Cells.Item(141, 11).Value = 28913.25
$ws.Cells.Item(141, 12).Value = 5199
$ws.Cells.Item(141, 13).Value = -23733.25
$ws.Cells.Item(141, 14).Value = -15559

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1460.9375
$ws.Cells.Item(2, 9).Value = 663.6
$ws.Cells.Item(2, 10).Value = 2789.8333
$ws.Cells.Item(2, 11).Value = 663.6
$ws.Cells.Item(2, 12).Value = 2789.8333
$ws.Cells.Item(2, 13).Value = -550.6
$ws.Cells.Item(2, 14).Value = -3015.8333
$ws.Cells.Item(61, 8).Value = 1094.7222
$ws.Cells.Item(61, 9).Value = 797.3570999999999
$ws.Cells.Item(61, 10).Value = 2135.5
$ws.Cells.Item(61, 11).Value = 797.3570999999999
$ws.Cells.Item(61, 12).Value = 2135.5
$ws.Cells.Item(61, 13).Value = -585.3570999999999
$ws.Cells.Item(61, 14).Value = -2559.5
$ws.Cells.Item(116, 8).Value = 1460.9375
$ws.Cells.Item(116, 9).Value = 663.6
$ws.Cells.Item(116, 10).Value = 2789.8333
$ws.Cells.Item(116, 11).Value = 663.6
$ws.Cells.Item(116, 12).Value = 2789.8333
$ws.Cells.Item(116, 13).Value = 1630.4
$ws.Cells.Item(116, 14).Value = -7377.8333
$ws.Cells.Item(132, 8).Value = 2072.1702
$ws.Cells.Item(132, 9).Value = 1776.0476
$ws.Cells.Item(132, 11).Value = 5328.142800000001
$ws.Cells.Item(132, 13).Value = -2798.142800000001
$ws.Cells.Item(136, 8).Value = 1094.7222
$ws.Cells.Item(136, 9).Value = 797.3570999999999
$ws.Cells.Item(136, 10).Value = 2135.5
$ws.Cells.Item(136, 11).Value = 2392.0713
$ws.Cells.Item(136, 12).Value = 6406.5
$ws.Cells.Item(136, 13).Value = 157.9287000000004
$ws.Cells.Item(136, 14).Value = -11506.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1460.9375
$ws.Cells.Item(3, 9).Value = 663.6
$ws.Cells.Item(3, 10).Value = 2789.8333
$ws.Cells.Item(3, 11).Value = 663.6
$ws.Cells.Item(3, 12).Value = 2789.8333
$ws.Cells.Item(3, 13).Value = -549.6
$ws.Cells.Item(3, 14).Value = -3017.8333
$ws.Cells.Item(134, 8).Value = 4633.606
$ws.Cells.Item(134, 9).Value = 1615.5714
$ws.Cells.Item(134, 11).Value = 4846.7142
$ws.Cells.Item(134, 13).Value = -2311.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 6629707.5
$ws.Cells.Item(6, 9).Value = 9939560
$ws.Cells.Item(6, 11).Value = 9939560
$ws.Cells.Item(6, 13).Value = -9939447
$ws.Cells.Item(16, 8).Value = 83334630
$ws.Cells.Item(16, 9).Value = 142858320
$ws.Cells.Item(16, 10).Value = 1457
$ws.Cells.Item(16, 11).Value = 142858320
$ws.Cells.Item(16, 12).Value = 1457
$ws.Cells.Item(16, 13).Value = -142858033
$ws.Cells.Item(16, 14).Value = -2031
$ws.Cells.Item(31, 8).Value = 1062.3928
$ws.Cells.Item(31, 9).Value = 764.2632
$ws.Cells.Item(31, 10).Value = 1691.7778
$ws.Cells.Item(31, 11).Value = 764.2632
$ws.Cells.Item(31, 12).Value = 1691.7778
$ws.Cells.Item(31, 13).Value = -469.2632
$ws.Cells.Item(31, 14).Value = -2281.7778
$ws.Cells.Item(34, 8).Value = 1062.3928
$ws.Cells.Item(34, 9).Value = 764.2632
$ws.Cells.Item(34, 10).Value = 1691.7778
$ws.Cells.Item(34, 11).Value = 764.2632
$ws.Cells.Item(34, 12).Value = 1691.7778
$ws.Cells.Item(34, 13).Value = -562.2632
$ws.Cells.Item(34, 14).Value = -2095.7778
$ws.Cells.Item(56, 8).Value = 0
$ws.Cells.Item(56, 9).Value = 0
$ws.Cells.Item(56, 10).Value = 0
$ws.Cells.Item(56, 11).Value = 0
$ws.Cells.Item(56, 12).Value = 0
$ws.Cells.Item(56, 13).ClearContents()
$ws.Cells.Item(56, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 83334630
$ws.Cells.Item(113, 9).Value = 142858320
$ws.Cells.Item(113, 10).Value = 1457
$ws.Cells.Item(113, 11).Value = 142858320
$ws.Cells.Item(113, 12).Value = 1457
$ws.Cells.Item(113, 13).Value = -142856150
$ws.Cells.Item(113, 14).Value = -5797

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(39, 8).Value = 0
$ws.Cells.Item(39, 9).Value = 0
$ws.Cells.Item(39, 10).Value = 0
$ws.Cells.Item(39, 11).Value = 0
$ws.Cells.Item(39, 12).Value = 0
$ws.Cells.Item(39, 13).ClearContents()
$ws.Cells.Item(39, 14).ClearContents()
$ws.Cells.Item(68, 8).Value = 2278.5881
$ws.Cells.Item(68, 9).Value = 849.4286
$ws.Cells.Item(68, 11).Value = 2548.2858
$ws.Cells.Item(68, 13).Value = -1737.2858
$ws.Cells.Item(71, 8).Value = 2278.5881
$ws.Cells.Item(71, 9).Value = 849.4286
$ws.Cells.Item(71, 11).Value = 7644.8574
$ws.Cells.Item(71, 13).Value = -3588.8574
$ws.Cells.Item(107, 8).Value = 3425.8057
$ws.Cells.Item(107, 9).Value = 555.7083
$ws.Cells.Item(107, 11).Value = 1667.1249
$ws.Cells.Item(107, 13).Value = 252.8751
$ws.Cells.Item(137, 8).Value = 10341.1
$ws.Cells.Item(137, 9).Value = 2387.5
$ws.Cells.Item(137, 11).Value = 7162.5
$ws.Cells.Item(137, 13).Value = -2062.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(133, 8).Value = 38874.25
$ws.Cells.Item(133, 10).Value = 38874.25
$ws.Cells.Item(133, 12).Value = 38874.25
$ws.Cells.Item(133, 14).Value = -48994.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 2590.7693
$ws.Cells.Item(40, 9).Value = 2334.0908
$ws.Cells.Item(40, 11).Value = 2334.0908
$ws.Cells.Item(40, 13).Value = -2198.0908
$ws.Cells.Item(61, 8).Value = 1992.6
$ws.Cells.Item(61, 9).Value = 1374.7858
$ws.Cells.Item(61, 11).Value = 1374.7858
$ws.Cells.Item(61, 13).Value = -1172.7858
$ws.Cells.Item(68, 8).Value = 1317.8125
$ws.Cells.Item(68, 9).Value = 1058.7
$ws.Cells.Item(68, 11).Value = 1058.7
$ws.Cells.Item(68, 13).Value = -309.7
$ws.Cells.Item(71, 8).Value = 1317.8125
$ws.Cells.Item(71, 9).Value = 1058.7
$ws.Cells.Item(71, 11).Value = 5293.5
$ws.Cells.Item(71, 13).Value = -1549.5
$ws.Cells.Item(82, 8).Value = 2396
$ws.Cells.Item(82, 9).Value = 2348
$ws.Cells.Item(82, 10).Value = 2540
$ws.Cells.Item(82, 11).Value = 2348
$ws.Cells.Item(82, 12).Value = 2540
$ws.Cells.Item(82, 13).Value = -1987
$ws.Cells.Item(82, 14).Value = -3262
$ws.Cells.Item(85, 8).Value = 2396
$ws.Cells.Item(85, 9).Value = 2348
$ws.Cells.Item(85, 10).Value = 2540
$ws.Cells.Item(85, 11).Value = 2348
$ws.Cells.Item(85, 12).Value = 2540
$ws.Cells.Item(85, 13).Value = -1100
$ws.Cells.Item(85, 14).Value = -5036
$ws.Cells.Item(113, 8).Value = 1992.6
$ws.Cells.Item(113, 9).Value = 1374.7858
$ws.Cells.Item(113, 11).Value = 1374.7858
$ws.Cells.Item(113, 13).Value = 795.2141999999999
$ws.Cells.Item(136, 8).Value = 6063.6665
$ws.Cells.Item(136, 9).Value = 8594.462
$ws.Cells.Item(136, 10).Value = 1951.125
$ws.Cells.Item(136, 11).Value = 25783.386
$ws.Cells.Item(136, 12).Value = 5853.375
$ws.Cells.Item(136, 13).Value = -23233.386
$ws.Cells.Item(136, 14).Value = -10953.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 41671810
$ws.Cells.Item(62, 9).Value = 45458520
$ws.Cells.Item(62, 10).Value = 18000
$ws.Cells.Item(62, 11).Value = 45458520
$ws.Cells.Item(62, 12).Value = 18000
$ws.Cells.Item(62, 13).Value = -45457896
$ws.Cells.Item(62, 14).Value = -19248
$ws.Cells.Item(65, 8).Value = 41671810
$ws.Cells.Item(65, 9).Value = 45458520
$ws.Cells.Item(65, 10).Value = 18000
$ws.Cells.Item(65, 11).Value = 227292600
$ws.Cells.Item(65, 12).Value = 90000
$ws.Cells.Item(65, 13).Value = -227289480
$ws.Cells.Item(65, 14).Value = -96240
$ws.Cells.Item(80, 8).Value = 14920
$ws.Cells.Item(80, 10).Value = 14920
$ws.Cells.Item(80, 12).Value = 14920
$ws.Cells.Item(80, 14).Value = -16916
$ws.Cells.Item(81, 8).Value = 3119.2856
$ws.Cells.Item(81, 9).Value = 2816.8
$ws.Cells.Item(81, 10).Value = 3875.5
$ws.Cells.Item(81, 11).Value = 5633.6
$ws.Cells.Item(81, 12).Value = 7751
$ws.Cells.Item(81, 13).Value = -4572.6
$ws.Cells.Item(81, 14).Value = -9873
$ws.Cells.Item(83, 8).Value = 14920
$ws.Cells.Item(83, 10).Value = 14920
$ws.Cells.Item(83, 12).Value = 44760
$ws.Cells.Item(83, 14).Value = -54744
$ws.Cells.Item(84, 8).Value = 3119.2856
$ws.Cells.Item(84, 9).Value = 2816.8
$ws.Cells.Item(84, 10).Value = 3875.5
$ws.Cells.Item(84, 11).Value = 28168
$ws.Cells.Item(84, 12).Value = 38755
$ws.Cells.Item(84, 13).Value = -22864
$ws.Cells.Item(84, 14).Value = -49363
$ws.Cells.Item(124, 8).Value = 52439.668
$ws.Cells.Item(124, 10).Value = 52439.668
$ws.Cells.Item(124, 12).Value = 52439.668
$ws.Cells.Item(124, 14).Value = -62259.668
$ws.Cells.Item(132, 8).Value = 4221
$ws.Cells.Item(132, 9).Value = 5226.1177
$ws.Cells.Item(132, 11).Value = 15678.3531
$ws.Cells.Item(132, 13).Value = -13148.3531
